$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "Volta a la Comunitat Valenciana"
$ws.Range("B13").Value = "Volta a la Comunitat Valenciana"

$ws.Range("B19").Select()
